# Fix old holder file
# Refresh the WhaleWatching holders sheet: ranks 12-14 get new
# balances/addresses, rank 23's balance ticks up slightly, the
# "reshuffled" highlight moves from ranks 24/43 to rank 14, and the
# Totals row is recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must remain plain TEXT (Excel would
# otherwise auto-detect large comma-grouped numbers / bare percentages
# and silently convert the cell to a numeric type with a new number
# format). We briefly force a Text format, set the value, then restore
# the original look of the cell by pasting the formatting from a
# neighboring cell that already carries the desired style, so the
# cell's style index ends up unchanged.
function Set-TextValue($targetAddr, $value, $formatSourceAddr) {
    $cell = $ws.Range($targetAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $ws.Range($formatSourceAddr).Copy()
    $cell.PasteSpecial(-4122)
    $ws.Application.CutCopyMode = $false
}

# ---------------------------------------------------------------------
# Step 1: Copy the existing "reshuffled row" highlight formatting (from
# row 44, which currently owns it) onto row 15 before that formatting
# gets removed from rows 25 / 44 further down.
# ---------------------------------------------------------------------
$ws.Range("A44:E44").Copy()
$ws.Range("A15:E15").PasteSpecial(-4122)

$ws.Range("F44:G44").Copy()
$ws.Range("F15:G15").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# Step 2: Update holder data for rows 13-15 (ranks 12-14).
# ---------------------------------------------------------------------

# Row 13 (rank 12)
$ws.Range("B13").Value = "0x7167d70b2e6b167633356254bf22f5d32aedcd5b"
Set-TextValue "C13" "405,964,100,208,639" "C14"
$ws.Range("D13").Value = "406 Trillion"
$ws.Range("E13").Value = "0.4060% "

# Row 14 (rank 13)
$ws.Range("B14").Value = "0x7395cb62e4405b6c1174c2329f444af6ee7bdfd9"
Set-TextValue "C14" "362,521,477,411,729" "C13"
$ws.Range("D14").Value = "363 Trillion"
$ws.Range("E14").Value = "0.3625% "

# Row 15 (rank 14) - now the highlighted "changed" row
# (Use row 44's still-highlighted C/G cells as the format reference so
# the right style - shaded + right-aligned [+ bold for G] - is kept.)
$ws.Range("B15").Value = "0x69fe97ce030074b37cbaf3ee46e9f68ca8712099"
Set-TextValue "C15" "360,003,463,717,642" "C44"
$ws.Range("D15").Value = "360 Trillion"
$ws.Range("E15").Value = "0.3600% "
$ws.Range("F15").Value = "47 Trillion Less"
Set-TextValue "G15" "46,742,660,549,411" "G44"

# ---------------------------------------------------------------------
# Step 3: Update row 24's amount (rank 23).
# ---------------------------------------------------------------------
Set-TextValue "C24" "170,559,005,253,847" "C23"

# ---------------------------------------------------------------------
# Step 4: Row 25 (rank 24) is no longer the reshuffled/highlighted row,
# so restore its formatting to match the surrounding normal rows and
# clear out the "More/Less" columns.
# ---------------------------------------------------------------------
$ws.Range("A24:E24").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)

$ws.Range("F24:G24").Copy()
$ws.Range("F25:G25").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

$ws.Range("F25").ClearContents()
$ws.Range("G25").ClearContents()

# ---------------------------------------------------------------------
# Step 5: Row 44 (rank 43) is also no longer highlighted; restore its
# formatting and clear the "More/Less" columns.
# ---------------------------------------------------------------------
$ws.Range("A43:E43").Copy()
$ws.Range("A44:E44").PasteSpecial(-4122)

$ws.Range("F43:G43").Copy()
$ws.Range("F44:G44").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

$ws.Range("F44").ClearContents()
$ws.Range("G44").ClearContents()

# ---------------------------------------------------------------------
# Step 6: Update the Totals row (53) to reflect the new aggregate
# numbers.
# ---------------------------------------------------------------------
Set-TextValue "C53" "12,035,571,222,282,294" "C24"
$ws.Range("D53").Value = "12036 Trillion"
Set-TextValue "E53" "12.04%" "B53"
$ws.Range("F53").Value = "-47 Trillion"
Set-TextValue "G53" "-46,742,660,549,411" "F53"
